$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.062.36"
$ws.Range("E2").Value = "  -1.49%  "
$ws.Range("D3").Value = "1.793.43"
$ws.Range("E3").Value = "  -1.87%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.32%  "
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.29"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.03"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("E10").Value = "  -2.17%  "
$ws.Range("E11").Value = "  -3.59%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0927"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.29%  "
$ws.Range("D13").Value = "2.051.18"
$ws.Range("E13").Value = "  -1.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.38"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.99%  "
$ws.Range("D15").Value = "1.790.51"
$ws.Range("E15").Value = "  -1.98%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.636"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.86%  "
$ws.Range("D17").Value = "34.078.75"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.22"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.84%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.62"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "253.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.10%  "
$ws.Range("E21").Value = "  -2.55%  "
$ws.Range("E22").Value = "  +0.37%  "
$ws.Range("E23").Value = "  -1.52%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.13"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.60"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.78%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.62"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.115"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.59%  "
$ws.Range("E29").Value = "  -2.76%  "
$ws.Range("E30").Value = "  +0.29%  "
$ws.Range("E31").Value = "  +0.34%  "
$ws.Range("E32").Value = "  -0.35%  "
$ws.Range("E33").Value = "  -0.79%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.63"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.85"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.00%  "
$ws.Range("D36").Value = "1.477.45"
$ws.Range("E36").Value = "  -7.82%  "
$ws.Range("E37").Value = "  -0.52%  "
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("E39").Value = "  -1.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "84.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.10%  "
$ws.Range("E41").Value = "  +0.46%  "
$ws.Range("E42").Value = "  -0.65%  "
$ws.Range("E43").Value = "  -2.68%  "
$ws.Range("E44").Value = "  -4.69%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0516"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.34%  "
$ws.Range("E46").Value = "  +0.45%  "
$ws.Range("D47").Value = "1.949.78"
$ws.Range("E47").Value = "  -1.67%  "
$ws.Range("E48").Value = "  +0.26%  "
$ws.Range("E49").Value = "  -1.48%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.85"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.59"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.50%  "
